# Update Generic Excel DataLayer (iPasXL)
# Re-populate the "Equipment" sample data: refresh row 2, replace row 3 with
# the data that used to live in row 4, and remove the old rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - new sampled values
$ws.Range("B2").Value = "DESC-7"
$ws.Range("C2").Value = "PT-9"
$ws.Range("D2").Value = "PDT-9"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2

# Row 3 - new sampled values
$ws.Range("A3").Value = "Equip-003"
$ws.Range("B3").Value = "DESC-5"
$ws.Range("C3").Value = "PT-5"
$ws.Range("D3").Value = "PDT-6"
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3

# Remove old rows 4 and 5 entirely, shrinking the table to 2 data rows
$ws.Rows("4:5").Delete()
